$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T45")

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.03459060192254842
$ws.Range("C2").Value = 2.120780142017438
$ws.Range("D2").Value = 17.47281255924218
$ws.Range("E2").Value = 4.180049348900343
$ws.Range("F2").Value = 4.278270417603791
$ws.Range("G2").Value = 22

# Row 3 (Q0)
$ws.Range("B3").Value = 0.9554726533143864
$ws.Range("C3").Value = 2.224672895214267
$ws.Range("D3").Value = 20.02789313565289
$ws.Range("E3").Value = 4.475253415802605
$ws.Range("F3").Value = 4.470327414903789
$ws.Range("G3").Value = 23

# Row 4 (Q1)
$ws.Range("B4").Value = -0.1311017249985545
$ws.Range("C4").Value = 0.9759568923264486
$ws.Range("D4").Value = 3.643887691671368
$ws.Range("E4").Value = 1.908896982990797
$ws.Range("F4").Value = 1.953862579171877
$ws.Range("G4").Value = 20
